$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Статистика по годам")

# Update column widths (C: 32 -> 28, E: 35 -> 31).
# Excel's ColumnWidth property (characters) gets converted internally to the
# stored "characters" width using the default font metrics, which adds an
# extra ~0.8333 fraction when fed a whole number (e.g. 28 -> 28.8333...).
# Feeding in (target - 11/12) reproduces the exact target width after
# Excel's internal round-trip conversion.
$ws.Columns.Item(3).ColumnWidth = 28 - 11/12 + 0.002
$ws.Columns.Item(5).ColumnWidth = 31 - 11/12 + 0.002

# Update header text
$ws.Range("C1").Value = "Средняя зарплата - Инженер"
$ws.Range("E1").Value = "Количество вакансий - Инженер"

# Update data values
$ws.Range("C2").Value = 68861
$ws.Range("E2").Value = 7
